# Applies the "Atualizado por script em 03-01-2024 20:45" update:
#  - Rows 31/32, 44/45, 58/60, 86/87: the betexplorer scraper re-ordered
#    pairs of same-kickoff-time matches, so the match data (columns F:V)
#    for each pair is swapped; columns A:E (index/pais/torneio/temporada/
#    data_partida) stay put.
#  - Four new match rows (94-97) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($rowA, $rowB)
    $rangeA = $ws.Range("F${rowA}:V${rowA}")
    $rangeB = $ws.Range("F${rowB}:V${rowB}")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-RowData 31 32
Swap-RowData 44 45
Swap-RowData 58 60
Swap-RowData 86 87

# Append the 4 new rows. Copy row 93's formatting (bold/bordered index
# column, date-formatted data_partida column) down onto 94:97 first, since
# the new rows don't exist yet and direct multi-cell writes to brand-new
# rows are unreliable -- write every new cell individually below.
$ws.Range("A93:V93").Copy()
$ws.Range("A94:V97").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @{
    94 = @(93, "israel", "ligat-ha-al", "2023-2024", 45294.75, "Hapoel Jerusalem", 3, "Maccabi Petah Tikva", 1, 2.63, "31/12/2024 19:12", 2.06, "03/01/2024 17:55", 3, "31/12/2024 19:12", 3.31, "03/01/2024 17:55", 2.78, "31/12/2024 19:12", 3.87, "03/01/2024 17:55", "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-jerusalem-maccabi-petah-tikva/K8DeG6SF/")
    95 = @(94, "israel", "ligat-ha-al", "2023-2024", 45294.78125, "Hapoel Tel Aviv", 3, "Sakhnin", 3, 2.31, "31/12/2024 19:42", 2.23, "03/01/2024 18:43", 3.38, "31/12/2024 19:42", 3.42, "03/01/2024 18:43", 2.92, "31/12/2024 19:42", 3.3, "03/01/2024 18:40", "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-tel-aviv-sakhnin/tIEiHnD9/")
    96 = @(95, "israel", "ligat-ha-al", "2023-2024", 45294.79166666666, "SC Ashdod", 1, "Maccabi Tel Aviv", 4, 9.51, "31/12/2024 19:12", 9.279999999999999, "03/01/2024 18:56", 5.57, "31/12/2024 19:12", 5.22, "03/01/2024 18:56", 1.24, "31/12/2024 19:12", 1.33, "03/01/2024 18:56", "https://www.betexplorer.com/football/israel/ligat-ha-al/sc-ashdod-maccabi-tel-aviv/WnX0fP5d/")
    97 = @(96, "israel", "ligat-ha-al", "2023-2024", 45294.80208333334, "Hapoel Petah Tikva", 2, "Maccabi Haifa", 2, 8.15, "31/12/2024 19:42", 10.3, "03/01/2024 19:10", 5.22, "31/12/2024 19:42", 5.71, "03/01/2024 19:10", 1.29, "31/12/2024 19:42", 1.28, "03/01/2024 19:09", "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-petah-tikva-maccabi-haifa/IqTde5jj/")
}

foreach ($r in 94..97) {
    $rowVals = $newRows[$r]
    for ($i = 0; $i -lt $rowVals.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value2 = $rowVals[$i]
    }
}
